# TARUGO ESPIGA DISMAY price list refresh:
# - bump the list's date stamp (A1) by one month
# - update the four bag prices (D32:D35) to the new values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 45436

$ws.Range("D32").Value = 13166.095
$ws.Range("D33").Value = 10457.641
$ws.Range("D34").Value = 10432.562
$ws.Range("D35").Value = 15009.349
